# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price (column D) values are stored as literal text (they keep trailing
# zeros like "609.72" / "1.00"), so numeric-looking replacements are typed
# with a leading apostrophe to force text entry, exactly as Excel would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''64.341.54'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '''3.141.88'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''609.72'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').Value = '''143.62'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''3.139.04'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').Value = '''0.530'
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').Value = '''0.151'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').Value = '''5.44'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '''0.477'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('E13').Value = '  +2.69%  '
$ws.Range('D14').Value = '''35.60'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = '''3.655.52'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '''64.243.88'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '''3.120.11'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').Value = '''6.89'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').Value = '''477.50'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = '''14.74'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').Value = '''0.725'
$ws.Range('E22').Value = '  +2.49%  '
$ws.Range('D23').Value = '''7.85'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').Value = '''13.70'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''85.40'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  +0.03%  '

# Rows 27/28 swapped rank order (RenderToken moved above PancakeSwap).
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '''8.59'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = '''2.79'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('E29').Value = '  +9.38%  '
$ws.Range('E30').Value = '  +3.58%  '
$ws.Range('D31').Value = '''2.09'
$ws.Range('E31').Value = '  -4.64%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = '''26.63'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('D34').Value = '''2.65'
$ws.Range('E34').Value = '  -3.11%  '
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = '''5.96'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '''52.39'
$ws.Range('E37').Value = '  -3.47%  '
$ws.Range('D38').Value = '''0.0₃0747'
$ws.Range('E38').Value = '  +5.17%  '
$ws.Range('D39').Value = '''454.62'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').Value = '''3.02'
$ws.Range('E40').Value = '  +5.17%  '
$ws.Range('D41').Value = '''0.0397'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('D43').Value = '''8.35'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '''2.879.12'
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('D45').Value = '''0.265'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').Value = '''2.27'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').Value = '''2.45'
$ws.Range('E47').Value = '  +6.14%  '
$ws.Range('D48').Value = '''26.59'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '''120.90'
$ws.Range('E51').Value = '  +2.19%  '
